# Insert a new weekly data row at row 170, shifting existing rows 170-274
# down to 171-275, then populate the new row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 170 (pushes rows 170..274 -> 171..275)
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record's data
$ws.Cells.Item(170, 1).Value = 3
$ws.Cells.Item(170, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = 44582
$ws.Cells.Item(170, 5).Value = 5
$ws.Cells.Item(170, 6).Value = 100112009
$ws.Cells.Item(170, 7).Value = "Acelga"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 230
$ws.Cells.Item(170, 11).Value = 2300
$ws.Cells.Item(170, 12).Value = 2500
$ws.Cells.Item(170, 13).Value = 2396
$ws.Cells.Item(170, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(170, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(170, 16).Value = 399
$ws.Cells.Item(170, 17).Value = 6
$ws.Cells.Item(170, 18).Value = "Hortaliza"
